$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the GitHub link text for a few rows ---
$githubLink = "https://github.com/KJGangarsha/ADAPT"

$ws.Range("F5").Value = $githubLink
$ws.Range("F12").Value = $githubLink
$ws.Range("F16").Value = $githubLink

# Apply a smaller (7pt), centered font to the cells that got the link as
# well as the two cells that only pick up the new style (no value change).
# (Ranges built from a multi-area address only reformat the first area in
# this engine, so touch each cell individually to be safe.)
foreach ($addr in @("F3", "F4", "F5", "F12", "F16")) {
    $cell = $ws.Range($addr)
    $cell.Font.Size = 7
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4108     # xlCenter
}

# --- Update the active selection on the sheet ---
$ws.Range("E16").Select() | Out-Null

# --- Page setup: paper size (A4) + orientation (portrait) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
